$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (so numeric-looking strings like "26000.0" / "17080824033178231" keep
# their exact textual representation instead of being coerced into a
# floating point number and losing the trailing ".0" / precision).
#
# We stage the text in a scratch cell (far outside the used range) with
# NumberFormat "@" applied there, then Copy/PasteSpecial the scratch cell
# onto the real target. PasteSpecial only moves the *value* here, so the
# destination cell keeps whatever style it already had (its column default,
# or an explicit style it was carrying before), instead of inheriting the
# scratch cell's "@" formatted style.
function Set-TextValue($cellRef, $value) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4104)
}

# ---- Row 2 updates ----
Set-TextValue "G2" "26000.0"
Set-TextValue "I2" "+10000"
$ws.Range("J2").Value = " on 2024-02-16 at 16:46:59"
$ws.Range("K2").Value = "Received 10000 from Shubhash Singh"

# ---- Row 3 updates ----
Set-TextValue "G3" "11000"
Set-TextValue "H3" "10000"
Set-TextValue "I3" "+10000"
$ws.Range("J3").Value = "2024-02-16 at 19:50:01"
$ws.Range("K3").Value = "+10000 Deposite on"

# ---- Row 4 (new row) ----
$ws.Range("A4").Value = "Shubhash Singh"
Set-TextValue "B4" "17080824033178231"
$ws.Range("C4").Value = "+91-9381560406"
$ws.Range("D4").Value = "subhashsingh2059@gmail.com"
$ws.Range("E4").Value = "Zapkto9@"
$ws.Range("F4").Value = "Saving"
Set-TextValue "G4" "120000"
Set-TextValue "H4" "120000"
$ws.Range("I4").Value = "+ 120000"
$ws.Range("J4").Value = "2024-02-16 at 16:50:04"
$ws.Range("K4").Value = "+120000 Initial Deposite on"

# Clean up the scratch cell so it doesn't leave stray data/dimension changes.
$ws.Range("ZZ1").Clear()
